$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.11550149516097719
$ws.Range("A2").Value = -0.0059999999736781717
$ws.Range("A3").Value = 0.027874155385164912
$ws.Range("A4").Value = -0.0079999999601074734
$ws.Range("A5").Value = -0.0029999999797682975
$ws.Range("A6").Value = -0.0019999999797732926
$ws.Range("A7").Value = -0.0099999999447679677
$ws.Range("A8").Value = -0.009999999943181237
$ws.Range("A9").Value = -0.0019999999758333331
$ws.Range("A10").Value = -0.0019999999743234298
$ws.Range("A11").Value = -0.0029999999700347502
$ws.Range("A12").Value = -0.0034999999675413562
$ws.Range("A13").Value = -0.0034999999663600789
$ws.Range("A14").Value = 0.0025885926896167533
$ws.Range("A15").Value = -0.00099999997641297966
$ws.Range("A16").Value = -0.001999999971661115
$ws.Range("A17").Value = -0.0019999999710345051
$ws.Range("A18").Value = -0.0039999999622590821
$ws.Range("A19").Value = -0.00399999998216094
$ws.Range("A20").Value = -0.0039999999808699727
$ws.Range("A21").Value = -0.0039999999806701325
$ws.Range("A22").Value = -0.0039999999805155895
$ws.Range("A23").Value = -0.0049999999733856271
$ws.Range("A24").Value = -0.01999999990641399
$ws.Range("A25").Value = -0.019999999905198074
$ws.Range("A26").Value = -0.0024999999774006909
$ws.Range("A27").Value = -0.0024999999770769499
$ws.Range("A28").Value = -0.0019999999780306865
$ws.Range("A29").Value = 0.070429069583656023
$ws.Range("A30").Value = -0.05999999972521719
$ws.Range("A31").Value = -0.0069999999512013744
$ws.Range("A32").Value = -0.009999999937994275
$ws.Range("A33").Value = -0.0039999999636890493
